$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the source diff.
# Column D values are forced to Text format to prevent Excel/COM
# from auto-converting numeric-looking strings (e.g. "21.18", "0.00000000110")
# into floating point numbers, then the style is reset back to the default
# "Normal" style so no unintended style index is introduced on the cell.
$updates = @(
    @("D2", "26.121.83"),
    @("E2", "  -0.69%  "),
    @("D3", "1.667.32"),
    @("E3", "  -1.39%  "),
    @("E4", "  -0.53%  "),
    @("D5", "210.45"),
    @("E5", "  -3.38%  "),
    @("D6", "0.5251"),
    @("E7", "  -0.47%  "),
    @("D8", "0.2625"),
    @("E8", "  -3.91%  "),
    @("D9", "0.06289"),
    @("E9", "  -2.42%  "),
    @("D10", "21.18"),
    @("E10", "  -2.15%  "),
    @("D11", "0.07549"),
    @("D12", "1.681.69"),
    @("E12", "  -0.58%  "),
    @("D13", "4.437"),
    @("E13", "  -2.20%  "),
    @("D14", "0.5531"),
    @("E14", "  -4.46%  "),
    @("D15", "66.74"),
    @("E15", "  -0.54%  "),
    @("D16", "0.000007918"),
    @("E16", "  -5.51%  "),
    @("D17", "26.155.05"),
    @("E17", "  -0.72%  "),
    @("E18", "  -0.49%  "),
    @("D19", "4.735"),
    @("E19", "  -3.50%  "),
    @("D20", "186.21"),
    @("E20", "  -2.10%  "),
    @("D21", "10.33"),
    @("E21", "  -4.91%  "),
    @("D22", "6.167"),
    @("E22", "  -1.53%  "),
    @("E23", "  -0.55%  "),
    @("D24", "149.91"),
    @("E24", "  +0.55%  "),
    @("D25", "0.1248"),
    @("E25", "  -3.03%  "),
    @("D26", "7.490"),
    @("E26", "  -4.80%  "),
    @("D27", "15.92"),
    @("E27", "  +0.16%  "),
    @("D28", "0.06248"),
    @("E28", "  -0.43%  "),
    @("D29", "1.353"),
    @("E29", "  -1.43%  "),
    @("D30", "1.278"),
    @("E30", "  -3.64%  "),
    @("D31", "3.513"),
    @("E31", "  -2.41%  "),
    @("D32", "3.412"),
    @("E32", "  -4.74%  "),
    @("D33", "1.627"),
    @("E33", "  -2.87%  "),
    @("D34", "0.9970"),
    @("E34", "  -3.43%  "),
    @("D35", "0.6038"),
    @("E35", "  -1.86%  "),
    @("D36", "2.414"),
    @("E36", "  +0.01%  "),
    @("D37", "2.733"),
    @("E37", "  -1.18%  "),
    @("D38", "6.124"),
    @("E38", "  +0.19%  "),
    @("D39", "1.105.50"),
    @("E39", "  -0.36%  "),
    @("D40", "0.01613"),
    @("E40", "  -2.48%  "),
    @("D41", "0.8706"),
    @("E41", "  -1.25%  "),
    @("D42", "1.004"),
    @("D43", "99.98"),
    @("E43", "  -1.24%  "),
    @("D44", "1.821.58"),
    @("E44", "  -1.13%  "),
    @("B45", "BabyDogeCoin"),
    @("C45", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D45", "0.00000000110"),
    @("E45", "  -0.05%  "),
    @("B46", "Aave"),
    @("C46", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("D46", "55.37"),
    @("E46", "  -3.90%  "),
    @("B47", "Frax"),
    @("C47", "https://coinranking.com/coin/KfWtaeV1W+frax-frax"),
    @("D47", "1.004"),
    @("E47", "  -0.02%  "),
    @("B48", "EnergySwap"),
    @("C48", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D48", "8.056"),
    @("E48", "  -1.18%  "),
    @("B49", "Cronos"),
    @("C49", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"),
    @("D49", "0.05234"),
    @("E49", "  -0.88%  "),
    @("B50", "Mantle"),
    @("C50", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"),
    @("D50", "0.4246"),
    @("E50", "  -1.21%  "),
    @("B51", "Aptos"),
    @("C51", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"),
    @("D51", "5.964"),
    @("E51", "  -1.32%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $col = $addr.Substring(0,1)
    $range = $ws.Range($addr)
    if ($col -eq "D") {
        $range.NumberFormat = "@"
        $range.Value = $val
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}

Write-Host "Done applying" $updates.Length "cell updates."